# Added some error checking
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "MOCK_DATA (1)"

# Update the active selection on the sheet
$ws.Range("B30").Select()
